$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Cats Slot Game for Free | IGT
#    Online Slots"). The new paragraph must NOT inherit the Heading1
#    style, and must reproduce the same run layout used elsewhere in
#    the document: an empty leading run, a bold "Meta description"
#    run, and a plain run with the rest of the sentence.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($p1.Range.End, $p1.Range.End)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the wild savannah with Cats, a 30-payline slot game by IGT. Play for free and win big with majestic big cats and exciting features.</w:t></w:r></w:p>' +
    '<w:p/>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($metaXml)

# InsertXML needed a trailing empty paragraph to force the split away
# from the following "Cats Rule in This Fun and Wild Slot Game"
# paragraph; remove that now-stray empty paragraph (paragraph 3).
$stray = $d.Paragraphs.Item(3)
$stray.Range.Delete()

# ---------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Cats Slot Game for Free | IGT
#    Online Slots" paragraph that used to sit near the end of the
#    document (right before the italic meta-description paragraph).
#    Locate it by content (searching from the end, since paragraph 1
#    legitimately keeps the same title text) rather than assuming a
#    fixed index.
# ---------------------------------------------------------------------
$dupIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r") -eq "Play Cats Slot Game for Free | IGT Online Slots") {
        $dupIndex = $i
        break
    }
}
if ($dupIndex -ge 2) {
    $d.Paragraphs.Item($dupIndex).Range.Delete()
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final italic paragraph (previously the
#    meta-description sentence) with the new image-prompt text. Build
#    a fresh Range from explicit Start/End (rather than reusing the
#    live Paragraph.Range, and rather than Find.Execute's replacement
#    path which mangles straight quotes into smart quotes) so the
#    paragraph's run formatting (italic) is kept and the text - with
#    its literal straight double-quotes - lands exactly as authored.
# ---------------------------------------------------------------------
$oldSentence = "Experience the wild savannah with Cats, a 30-payline slot game by IGT. Play for free and win big with majestic big cats and exciting features."
$descIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r") -eq $oldSentence) {
        $descIndex = $i
        break
    }
}
if ($descIndex -eq -1) {
    $descIndex = $d.Paragraphs.Count
}

$targetP = $d.Paragraphs.Item($descIndex)
$targetRange = $d.Range($targetP.Range.Start, $targetP.Range.End)
$targetRange.Text = 'Create a feature image for "Cats" slot game by IGT that showcases a happy Maya warrior wearing glasses in a cartoon style. The image should feature the warrior surrounded by some of the majestic big cats of the game, such as jaguars, lions, tigers, leopards, and mountain lions. The background should resemble an African savannah with a wild sunset landscape, reeds, and trees. The image should convey the theme of the game, which is a bold adventure into the animal world and the heart of wild nature, and also highlight the high winning potential of the game through the happy expression of the warrior and the playful nature of the cats.'
